# Add a "Version:" column (header + value) to row 1 of every template
# sheet, and move the active selection to A3 (first visible data row).
#
# Row 1 on each sheet is a hidden metadata row of the form:
#   A1 = "Data type:"   B1 = <sheet's data-type name>
# We extend it with:
#   C1 = "Version:"     D1 = 1   (numeric)

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("C1").Value = "Version:"
    $ws.Range("D1").Value = 1

    # Move the selection/active cell to A3, matching the updated template.
    $ws.Range("A3").Select()
}

# Range.Select() above also activates the sheet it's on, so the workbook's
# originally-active tab ("assembly", the first sheet) ends up being whatever
# sheet we touched last. Re-activate the original tab to restore it.
$wb.Worksheets.Item(1).Range("A3").Select()
